$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts Password..MobileNo right by one)
$ws.Columns("D:D").Insert()

# Populate the new "email" column header and value
$ws.Range("D1").Value = "email"
$ws.Range("D2").Value = "rojantest3@gmail.com"

# Widen the new column to fit its content
$ws.Range("D1").ColumnWidth = 18.45

# Move the active selection to D3 (matches saved selection state)
$ws.Range("D3").Select()
